$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 99; $r++) {
    if ($r -eq 36) { continue }

    $eCell = $ws.Cells.Item($r, 5)   # column E (剩余)
    $fCell = $ws.Cells.Item($r, 6)   # column F (开始时间)

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($eVal -eq 1) {
        $eCell.Value2 = 7
        $fCell.Value2 = $fVal + 7
    } else {
        $eCell.Value2 = $eVal - 1
    }
}
